$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 36: the date for the existing "Vejledning fra Anders" entry moves
# from 2020-05-20 (43971) to 2020-05-21 (43972)
$ws.Cells.Item(36, 3).Value = 43972

# New rows 37-40: four new time-registration entries on 2020-05-21
$ws.Cells.Item(37, 1).Value = "Lav SSD03 og review SSD10"
$ws.Cells.Item(37, 3).Value = 43972
$ws.Cells.Item(37, 4).Value = 0.354166666666667
$ws.Cells.Item(37, 5).Value = 0.458333333333333

$ws.Cells.Item(38, 1).Value = "Rettelse af SD01xx og DCD01xx"
$ws.Cells.Item(38, 3).Value = 43972
$ws.Cells.Item(38, 4).Value = 0.458333333333333
$ws.Cells.Item(38, 5).Value = 0.583333333333333

$ws.Cells.Item(39, 1).Value = "Lav DD03 og DD10"
$ws.Cells.Item(39, 3).Value = 43972
$ws.Cells.Item(39, 4).Value = 0.583333333333333
$ws.Cells.Item(39, 5).Value = 0.635416666666667

$ws.Cells.Item(40, 1).Value = "Implementering af designet fra mockup"
$ws.Cells.Item(40, 3).Value = 43972
$ws.Cells.Item(40, 4).Value = 0.635416666666667
$ws.Cells.Item(40, 5).Value = 0.6875

# Update the active selection to mirror the author's last position (E41)
$ws.Activate()
$ws.Range("E41").Select()
